$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OOCL")

# Copy formatting only from the last data row (28) into the new row (29)
$ws.Range("A28:D28").Copy() | Out-Null
$ws.Range("A29:D29").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# New port mapping row: Sydney / AUSYD
$ws.Range("A29").Value = "AUSYD"
$ws.Range("B29").Value = "Sydney, ,New South Wales,Australia,AUSYD"
$ws.Range("C29").Value = 738872886232233
$ws.Range("D29").Value = "Sydney"

# Re-sort the data range (A2:D29) ascending by column A (Port Code)
$rng = $ws.Range("A2:D29")
$ws.Sort.SortFields.Clear() | Out-Null
$ws.Sort.SortFields.Add($ws.Range("A2:A29")) | Out-Null
$ws.Sort.SetRange($rng) | Out-Null
$ws.Sort.Header = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlNo
$ws.Sort.Apply() | Out-Null

# Extend the AutoFilter to cover the new row
$ws.AutoFilterMode = $false
$ws.Range("A1:B29").AutoFilter(1) | Out-Null

# Keep the _FilterDatabase defined name in sync with the new filter range
$wb.Names.Item("OOCL!_FilterDatabase").RefersTo = "=OOCL!`$A`$1:`$B`$29"

# Update the saved selection
$ws.Range("C11").Select() | Out-Null
